$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (sheet shrinks from 199 to 184 rows)
$ws.Range("A185:A199").EntireRow.Delete() | Out-Null

# Write every cell in column A (rows 1-184) with the updated content
$ws.Cells.Item(1, 1).Value = 'Closures/remote learning and relocationsGoogle Tag Manager (noscript)End Google Tag Manager (noscript)You may be trying to access this site from a secured browser on the server. Please enable scripts and reload this page.Start of Site headerSkip to contentSearchSearchEnd of Site headerHome>The Department>Program directory>Emergencies and Natural Disasters>Closures/remote learning and relocationsClosures/remote learning and relocations- Example page content area using current sharepoint wrapperPage Content'
$ws.Cells.Item(2, 1).Value = 'li: This page lists early childhood services, schools and TAFEs that are currently closed.'
$ws.Cells.Item(3, 1).Value = '*Please note these closures are NOT all related to the Covid-19 pandemic.'
$ws.Cells.Item(4, 1).Value = 'li: See'
$ws.Cells.Item(5, 1).Value = 'Coronavirus advice'
$ws.Cells.Item(6, 1).Value = 'li: for the latest advice on the outbreak of COVID-19 (coronavirus).'
$ws.Cells.Item(7, 1).Value = 'On this pageSchool and early childhood service, TAFE closures and relocationsBus service cancellations or alterationsSchool and early childhood service, TAFE closures and relocations for Tuesday 28 July 2020South-Eastern Victoria RegionEarly childhood services'
$ws.Cells.Item(8, 1).Value = 'li: The Department has been advised of the following'
$ws.Cells.Item(9, 1).Value = 'early childhood service'
$ws.Cells.Item(10, 1).Value = 'li: closures:'
$ws.Cells.Item(11, 1).Value = 'li: Ada Mary A''beckett Children''s Centre Inc PORT MELBOURNE'
$ws.Cells.Item(12, 1).Value = 'li: Aqua Energy Creche SALE'
$ws.Cells.Item(13, 1).Value = 'li: Berwick Neighbourhood Centre - Marriott Waters LYNDHURST'
$ws.Cells.Item(14, 1).Value = 'li: Berwick Neighbourhood Centre (Timbarra Playroom) BERWICK'
$ws.Cells.Item(15, 1).Value = 'li: Brighton Grammar OSHClub BRIGHTON'
$ws.Cells.Item(16, 1).Value = 'li: Camp Australia - Geelong Grammar School Toorak Campus OSHC TOORAK'
$ws.Cells.Item(17, 1).Value = 'li: Camp Australia - Haileybury Edrington OSHC BERWICK'
$ws.Cells.Item(18, 1).Value = 'li: Camp Australia - Haileybury Newlands OSHC KEYSBOROUGH'
$ws.Cells.Item(19, 1).Value = 'li: Camp Australia - Loreto College Mandeville Hall OSHC TOORAK'
$ws.Cells.Item(20, 1).Value = 'li: Camp Australia - Mentone Grammar OSHC MENTONE'
$ws.Cells.Item(21, 1).Value = 'li: Camp Australia - Peninsula Grammar OSHC MOUNT ELIZA'
$ws.Cells.Item(22, 1).Value = 'li: Camp Australia - Sacred Heart Parish School - Sandringham OSHC  SANDRINGHAM'
$ws.Cells.Item(23, 1).Value = 'li: Camp Australia - St Kevins College Toorak OSHC TOORAK'
$ws.Cells.Item(24, 1).Value = 'li: Camp Australia - Toorak College OSHC MOUNT ELIZA'
$ws.Cells.Item(25, 1).Value = 'li: Carnegie Occasional Care Centre CARNEGIE'
$ws.Cells.Item(26, 1).Value = 'li: Castlefield Community Centre HAMPTON'
$ws.Cells.Item(27, 1).Value = 'li: Clark Street Children''s Centre PORT MELBOURNE'
$ws.Cells.Item(28, 1).Value = 'li: Cranbourne Community House CRANBOURNE'
$ws.Cells.Item(29, 1).Value = 'li: Elwood Children''s Centre ELWOOD'
$ws.Cells.Item(30, 1).Value = 'li: Fernwood Fitness Clayton CLAYTON'
$ws.Cells.Item(31, 1).Value = 'li: Fernwood Fitness Narre Warren NARRE WARREN'
$ws.Cells.Item(32, 1).Value = 'li: GESAC Childcare BENTLEIGH EAST'
$ws.Cells.Item(33, 1).Value = 'li: Glen Iris Creche GLEN IRIS'
$ws.Cells.Item(34, 1).Value = 'li: Goodlife Chelsea Heights CHELSEA HEIGHTS'
$ws.Cells.Item(35, 1).Value = 'li: Goodlife Fountain Gate NARRE WARREN'
$ws.Cells.Item(36, 1).Value = 'li: Goodlife Karingal KARINGAL'
$ws.Cells.Item(37, 1).Value = 'li: Hallam Community Centre Inc HALLAM'
$ws.Cells.Item(38, 1).Value = 'li: Hope Frankston Heights FRANKSTON'
$ws.Cells.Item(39, 1).Value = 'li: Melbourne Grammar School, Grimwade House CAULFIELD'
$ws.Cells.Item(40, 1).Value = 'li: Merkaz Bentleigh BENTLEIGH'
$ws.Cells.Item(41, 1).Value = 'li: Only About Children Cheltenham CHELTENHAM'
$ws.Cells.Item(42, 1).Value = 'li: Overport Primary School OSHC - Extend FRANKSTON'
$ws.Cells.Item(43, 1).Value = 'li: Paisley Park Early Learning Centre Cranbourne CRANBOURNE'
$ws.Cells.Item(44, 1).Value = 'li: Pavillion - Frankston & District Netball Association Inc FRANKSTON'
$ws.Cells.Item(45, 1).Value = 'li: Peninsula Aquatic Recreation Centre FRANKSTON'
$ws.Cells.Item(46, 1).Value = 'li: Play Zone - Highett HIGHETT'
$ws.Cells.Item(47, 1).Value = 'li: Playtime Warragul WARRAGUL'
$ws.Cells.Item(48, 1).Value = 'li: Sacre Coeur Oshclub GLEN IRIS'
$ws.Cells.Item(49, 1).Value = 'li: Smaland Springvale Ikea SPRINGVALE'
$ws.Cells.Item(50, 1).Value = 'li: Somerville Recreation Centre Child Care SOMERVILLE'
$ws.Cells.Item(51, 1).Value = 'li: Toorak Primary School OSHClub TOORAK'
$ws.Cells.Item(52, 1).Value = 'li: Upper Beaconsfield Community Early Learning Centre BEACONSFIELD UPPER'
$ws.Cells.Item(53, 1).Value = 'Schools closedTheDepartment hasbeen advised of the followingschool closures:'
$ws.Cells.Item(54, 1).Value = 'li: Cheltenham Secondary College, CHELTENHAM'
$ws.Cells.Item(55, 1).Value = 'li: Dandenong High School, DANDENONG'
$ws.Cells.Item(56, 1).Value = 'li: Fountain Gate Secondary College, NARRE WARREN'
$ws.Cells.Item(57, 1).Value = 'li: Gleneagles Secondary College, ENDEAVOUR HILLS'
$ws.Cells.Item(58, 1).Value = 'li: Grayling Primary School, CLYDE NORTH'
$ws.Cells.Item(59, 1).Value = 'li: Kingston Heath Primary School, CHELTENHAM'
$ws.Cells.Item(60, 1).Value = 'li: Leibler Yavneh College, ELSTERNWICK'
$ws.Cells.Item(61, 1).Value = 'li: Overport Primary School, FRANKSTON'
$ws.Cells.Item(62, 1).Value = 'li: Parkdale Primary School, PARKDALE'
$ws.Cells.Item(63, 1).Value = 'li: Thomas Mitchell Primary School, ENDEAVOUR HILLS'
$ws.Cells.Item(64, 1).Value = 'li: Westall Primary School, CLAYTON SOUTH'
$ws.Cells.Item(65, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood services'
$ws.Cells.Item(66, 1).Value = 'li: The Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(67, 1).Value = 'li: Aqualink Box Hill Creche BOX HILL'
$ws.Cells.Item(68, 1).Value = 'li: Aqualink Nunawading Creche FOREST HILL'
$ws.Cells.Item(69, 1).Value = 'li: Camp Australia - Camberwell Boys Grammar Junior School OSHC CANTERBURY'
$ws.Cells.Item(70, 1).Value = 'li: Camp Australia - Mountain Gate Primary School OSHC FERNTREE GULLY'
$ws.Cells.Item(71, 1).Value = 'li: Camp Australia - Scoresby Primary School OSHC SCORESBY'
$ws.Cells.Item(72, 1).Value = 'li: Camp Australia - St Clement of Rome School OSHC BULLEEN'
$ws.Cells.Item(73, 1).Value = 'li: Camp Australia - St Timothy''s School Vermont OSHC VERMONT'
$ws.Cells.Item(74, 1).Value = 'li: Camp Australia - Strathcona Baptist Girls Junior School OSHC CANTERBURY'
$ws.Cells.Item(75, 1).Value = 'li: Camp Australia - Templestowe Heights Primary School OSHC TEMPLESTOWE LOWER'
$ws.Cells.Item(76, 1).Value = 'li: Carey Donvale OSHClub DONVALE'
$ws.Cells.Item(77, 1).Value = 'li: Clever Kids Childcare - Ashburton ASHBURTON'
$ws.Cells.Item(78, 1).Value = 'li: Fitness First Doncaster (Playzone) DONCASTER'
$ws.Cells.Item(79, 1).Value = 'li: Flamingo Community Group WANTIRNA SOUTH'
$ws.Cells.Item(80, 1).Value = 'li: MakerDojo HAWTHORN'
$ws.Cells.Item(81, 1).Value = 'li: St Andrews Christian College Outside School Hours Care WANTIRNA SOUTH'
$ws.Cells.Item(82, 1).Value = 'li: Vermont Primary School Kindergarten VERMONT'
$ws.Cells.Item(83, 1).Value = 'li: Vermont Primary School Outside School Hours Child Care Service VERMONT'
$ws.Cells.Item(84, 1).Value = 'Schools closedThe Department hasbeen advisedof the followingschool closures:Ashwood High School, BURWOODAuburn High School, HAWTHORN EASTGlen Iris Primary School, GLEN IRISKerrimuir Primary School, BOX HILL NORTH'
$ws.Cells.Item(85, 1).Value = 'li: Mooroolbark College, MOOROOLBARK'
$ws.Cells.Item(86, 1).Value = 'Mount Waverley Primary School, MOUNT WAVERLEYVermont Primary School, VERMONTWattle Park Primary School, BURWOODTAFE'
$ws.Cells.Item(87, 1).Value = 'li: The Department has'
$ws.Cells.Item(88, 1).Value = 'not'
$ws.Cells.Item(89, 1).Value = 'li: been advised of any TAFE closures.'
$ws.Cells.Item(90, 1).Value = 'North-Western Victoria RegionEarly childhood services'
$ws.Cells.Item(91, 1).Value = 'li: The Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(92, 1).Value = 'li: Big Childcare - Fitzroy PS OSHC FITZROY'
$ws.Cells.Item(93, 1).Value = 'li: Bright Stars Early Years Child Care Centre EPPING'
$ws.Cells.Item(94, 1).Value = 'li: Camp Australia - Ivanhoe Grammar OSHC IVANHOE'
$ws.Cells.Item(95, 1).Value = 'li: Camp Australia - Ivanhoe Grammar Plenty Valley Campus OSHC MERNDA'
$ws.Cells.Item(96, 1).Value = 'li: Camp Australia - St Joseph''s Primary School - Mernda OSHC MERNDA'
$ws.Cells.Item(97, 1).Value = 'li: Camp Australia - Thomastown Meadows Primary School OSHC THOMASTOWN'
$ws.Cells.Item(98, 1).Value = 'li: Collingwood College Afterschool Care and Vacation Care Program COLLINGWOOD'
$ws.Cells.Item(99, 1).Value = 'li: Diamond Creek Community Centre DIAMOND CREEK'
$ws.Cells.Item(100, 1).Value = 'li: Diamond Valley Sports and Fitness Centre GREENSBOROUGH'
$ws.Cells.Item(101, 1).Value = 'li: Eltham Leisure Centre ELTHAM'
$ws.Cells.Item(102, 1).Value = 'li: Gisborne Montessori School NEW GISBORNE'
$ws.Cells.Item(103, 1).Value = 'li: Kangaroo Ground Primary Combined OSHC KANGAROO GROUND'
$ws.Cells.Item(104, 1).Value = 'li: New Futures Broadmeadows BROADMEADOWS'
$ws.Cells.Item(105, 1).Value = 'li: New Futures Epping EPPING'
$ws.Cells.Item(106, 1).Value = 'li: Nino Early Learning Adventures - Bundoora BUNDOORA'
$ws.Cells.Item(107, 1).Value = 'li: Pender''s Grove Primary School Combined OSHC THORNBURY'
$ws.Cells.Item(108, 1).Value = 'li: Richmond West Afterschool Care and Vacation Care Program RICHMOND'
$ws.Cells.Item(109, 1).Value = 'li: St Bernard''s Out of School Hours Care COBURG EAST'
$ws.Cells.Item(110, 1).Value = 'li: Tullamarine Early Learning Centre TULLAMARINE'
$ws.Cells.Item(111, 1).Value = 'Schools'
$ws.Cells.Item(112, 1).Value = 'li: The Department has'
$ws.Cells.Item(113, 1).Value = 'li: been advised of the following school closures:'
$ws.Cells.Item(114, 1).Value = 'Al Siraat College, EPPINGCharles La Trobe, MACLEOD WESTDiamond Valley College, DIAMOND CREEKEpping Secondary College, EPPINGFitzroy High School, FITZROYFitzroy Primary School, FITZROYGisborne Secondary College, GISBORNEGladstone Park Secondary College, GLADSTONE PARKGreenvalePrimary School, GREENVALELalor Secondary College, LALORMelbourne Girls College, RICHMONDNewbury Primary School, CRAGIEBURNNorthern College of Arts and Technology (NCAT), PRESTONPascoe Vale Girls Secondary College, OAK PARKPenders Grove Primary School, THORNBURYPrinces Hill Secondary College, PRINCES HILLRoxburgh College, ROXBURGH PARKSacred Heart School, FITZROYThornbury High School, THORNBURYTAFE'
$ws.Cells.Item(115, 1).Value = 'li: The Department has'
$ws.Cells.Item(116, 1).Value = 'not'
$ws.Cells.Item(117, 1).Value = 'li: been advised of any TAFE closures.'
$ws.Cells.Item(118, 1).Value = 'South-Western Victoria RegionThe Department has not been advised of any school, early childhood service or TAFE closures, or buscancellations.Early childhood services'
$ws.Cells.Item(119, 1).Value = 'li: The Department has been advised of the followi'
$ws.Cells.Item(120, 1).Value = 'ng early childhood service closures:'
$ws.Cells.Item(121, 1).Value = 'li: Aerotots Activity Centre WERRIBEE'
$ws.Cells.Item(122, 1).Value = 'li: Altona Meadows Community Centre Occasional Care ALTONA MEADOWS'
$ws.Cells.Item(123, 1).Value = 'li: Aquapulse Creche HOPPERS CROSSING'
$ws.Cells.Item(124, 1).Value = 'li: Big Childcare - Manor Lakes P-12 College OSHC WYNDHAM VALE'
$ws.Cells.Item(125, 1).Value = 'li: Big Childcare – Sydenham/Hillside Sydenham Campus SYDENHAM'
$ws.Cells.Item(126, 1).Value = 'li: Blackwood Street Neighbourhood House YARRAVILLE'
$ws.Cells.Item(127, 1).Value = 'li: Bluewater Leisure Centre Creche COLAC'
$ws.Cells.Item(128, 1).Value = 'li: Camp Australia - Baden Powell P-9 College Derrimut Heath Campus OSHC HOPPERS CROSSING'
$ws.Cells.Item(129, 1).Value = 'li: Camp Australia - Footscray City Primary School OSHC FOOTSCRAY'
$ws.Cells.Item(130, 1).Value = 'li: Camp Australia - Haileybury City Campus OSHC WEST MELBOURNE'
$ws.Cells.Item(131, 1).Value = 'li: Camp Australia - Melton Christian College OSHC Melton South'
$ws.Cells.Item(132, 1).Value = 'li: Cana Catholic Primary OSHClub HILLSIDE'
$ws.Cells.Item(133, 1).Value = 'li: Carranballac Jamieson OSHClub POINT COOK'
$ws.Cells.Item(134, 1).Value = 'li: Coragulac & District Kindergarten CORAGULAC'
$ws.Cells.Item(135, 1).Value = 'li: Eagle Stadium WERRIBEE'
$ws.Cells.Item(136, 1).Value = 'li: Energy Force Fitness Creche DRYSDALE'
$ws.Cells.Item(137, 1).Value = 'li: Fernwood Fitness Sydenham SYDENHAM'
$ws.Cells.Item(138, 1).Value = 'li: Fun 4 All Occasional Care Center WERRIBEE'
$ws.Cells.Item(139, 1).Value = 'li: Future Kids Child Care-West Tarneit TARNEIT'
$ws.Cells.Item(140, 1).Value = 'li: Genesis Maidstone MAIDSTONE'
$ws.Cells.Item(141, 1).Value = 'li: Goodlife Essendon Child Minding ESSENDON'
$ws.Cells.Item(142, 1).Value = 'li: Goodlife Geelong BELMONT'
$ws.Cells.Item(143, 1).Value = 'li: Goodlife Point Cook POINT COOK'
$ws.Cells.Item(144, 1).Value = 'li: Goodlife Taylors Lakes TAYLORS LAKES'
$ws.Cells.Item(145, 1).Value = 'li: Happy Feet ELC MELTON WEST'
$ws.Cells.Item(146, 1).Value = 'li: Highpoint Kinder Haven MARIBYRNONG'
$ws.Cells.Item(147, 1).Value = 'li: Hopetoun Early Years Centre FLEMINGTON'
$ws.Cells.Item(148, 1).Value = 'li: Keilor Basketball Stadium Creche KEILOR PARK'
$ws.Cells.Item(149, 1).Value = 'li: Kensington Neighbourhood House Inc KENSINGTON'
$ws.Cells.Item(150, 1).Value = 'li: Kids Club Kensington Early Learning Centre KENSINGTON'
$ws.Cells.Item(151, 1).Value = 'li: Kids on Collins MELBOURNE'
$ws.Cells.Item(152, 1).Value = 'li: Lowther Hall OSHClub ESSENDON'
$ws.Cells.Item(153, 1).Value = 'li: Maribyrnong Aquatic Centre Occasional Child Care MARIBYRNONG'
$ws.Cells.Item(154, 1).Value = 'li: New Futures Braybrook BRAYBROOK'
$ws.Cells.Item(155, 1).Value = 'li: North Sunshine Kindergarten SUNSHINE'
$ws.Cells.Item(156, 1).Value = 'li: Shuter Street Occasional Care MOONEE PONDS'
$ws.Cells.Item(157, 1).Value = 'li: South Kingsville Community Centre SOUTH KINGSVILLE'
$ws.Cells.Item(158, 1).Value = 'li: St Monica''s OSHClub MOONEE PONDS'
$ws.Cells.Item(159, 1).Value = 'li: Story House Early Learning Keilor Downs KEILOR DOWNS'
$ws.Cells.Item(160, 1).Value = 'li: Sunshine Leisure Centre SUNSHINE'
$ws.Cells.Item(161, 1).Value = 'li: Willaura Primary School OSHC ARARAT'
$ws.Cells.Item(162, 1).Value = 'li: Winifred Nance Kindergarten COLAC'
$ws.Cells.Item(163, 1).Value = 'li: Woodlea Early Education AINTREE'
$ws.Cells.Item(164, 1).Value = 'li: Yarraville Community Centre YARRAVILLE'
$ws.Cells.Item(165, 1).Value = 'Schools closed'
$ws.Cells.Item(166, 1).Value = 'li: The Department has been advised of the following school closures:'
$ws.Cells.Item(167, 1).Value = 'Al Taqwa College, TRUGANINABraybrook College, BRAYBROOKBuckley Park College, ESSENDONCatholic Regional College, SYDENHAMClonard Secondary College, GEELONG WESTCopperfield College Delahey Campus, DELAHEYCofferfield College Sydenham Campus,DELAHEY'
$ws.Cells.Item(168, 1).Value = 'li: Footscray High School, FOOTSCRAY'
$ws.Cells.Item(169, 1).Value = 'Grovedale West Primary School, GROVEDALEHoppers Crossing Secondary College, HOPPERS CROSSINGLowther Hall Anglican Grammar School, ESSENDONKeilor Downs Secondary College, KEILOR DOWNSManor Lakes P-12 College, WYNDHAM VALEMary McKillop College, WERRIBEEMelton Secondary College, MELTONNorthcote High School, NORTHCOTEPoint Cook Senior Secondary College, POINT COOKRosamond Specialist School, BRAYBROOKSt Alban''s Secondary College, ST ALBANSSt Brigid''s College, HORSHAMSunshine Heights Primary School, SUNSHINESydenham Hillside Primary School, SYDENHAMTarneit Senior College, TARNEITTaylors Lakes Secondary College, TAYLORS LAKESThomas Carr College, TARNEITVictorian College of the Arts, ALBERT PARKWarracknabeal Secondary College, WARRACKNABEALTAFE'
$ws.Cells.Item(170, 1).Value = 'li: The Department has'
$ws.Cells.Item(171, 1).Value = 'not'
$ws.Cells.Item(172, 1).Value = 'li: been advised of any TAFE closures.'
$ws.Cells.Item(173, 1).Value = 'The Department hasnotbeen advised of any school, early childhood service or TAFE closures, or buscancellations.Bus service cancellations or alterations'
$ws.Cells.Item(174, 1).Value = 'li: For Term 2 2020, schools bus services will continue to be provided to support student travel to schools where needed.'
$ws.Cells.Item(175, 1).Value = 'li: With learning from home arrangements, local principals are authorised to cancel or suspend school buses where not required.'
$ws.Cells.Item(176, 1).Value = 'li: Families and students are encouraged to contact their school directly to enquire if their bus service is continuing to be provided.'
$ws.Cells.Item(177, 1).Value = 'li: Find more about bus services:'
$ws.Cells.Item(178, 1).Value = 'School bus services in Term 2 (for schools)School bus services in Term 2 (for parents)VicRoadsRegional roadsBack to top'
$ws.Cells.Item(179, 1).Value = 'li: Last Update: 28 July 2020'
$ws.Cells.Item(180, 1).Value = 'Website navigationFor parentsEarly childhoodSchoolsTAFE and trainingThe Department<li>                                 <a href="/about/educationstate">Education State                                                                  </a>                             </li>Support linksHelp in your languageAccessibilityPrivacyCopyrightContact linksContact usReport a website issue'
$ws.Cells.Item(181, 1).Value = 'li: State Government of Victoria, Australia © 2019'
$ws.Cells.Item(182, 1).Value = 'li: We respectfully acknowledge the Traditional Owners of country throughout Victoria and pay respect to the ongoing living cultures of First Peoples.'
$ws.Cells.Item(183, 1).Value = 'li: Our website uses a free tool to translate into other languages. This tool is a guide and may not be accurate. For more, see:'
$ws.Cells.Item(184, 1).Value = 'Information in your languageOld Search Code PlaceholdCludo Code for www.education.vic.gov.auHotjar Tracking Code for www.education.vic.gov.auSR-1181393'
